$d = $word.ActiveDocument
$apos = [char]0x2019

# ---------------------------------------------------------------------------
# Paragraph 1: remove the spell-check proofErr splits around ".bmp" (x2),
# "library" and the final ".bmp" by replacing the spanning text with itself.
# ---------------------------------------------------------------------------

$needle1 = "schermpje .bmp afbeeldingen te kunnen laten zien en later ook hiermee fonts te kunnen weergeven en zo dus tekst te kunnen laten zien. Hierv"
$d.Content.Find.Execute($needle1, $true, $false, $false, $false, $false, $true, 1, $false, $needle1, 1) | Out-Null

$needle2 = "van de .bmp afbeeldingen om verwerkt te worden door "
$d.Content.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, $needle2, 1) | Out-Null

$needle3 = " library. Daarna ga ik een script in python schrijven die een .bmp kan omtoveren tot zo" + $apos + "n formaat (waarschijnlijk een array in een of andere vorm). "
$d.Content.Find.Execute($needle3, $true, $false, $false, $false, $false, $true, 1, $false, $needle3, 1) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 2: rewrite the sentence about extending the library, then remove
# the proofErr splits around "library" and "bitmap" (x2).
# ---------------------------------------------------------------------------

$needle4 = "het maken van de library zelf, waarbij ik eerst afbeeldingen kan laten zien op het scherm"
$d.Content.Find.Execute($needle4, $true, $false, $false, $false, $false, $true, 1, $false, $needle4, 1) | Out-Null

$needle5 = " en vervolgens dit kan uitbreiden naar fonts en dus ook tekst."
$replace5 = " en hier bewerkingen op kan uitvoeren, zoals inverteren en croppen. Vervolgens kan ik dit uitbreiden naar fonts en dus tekst."
$d.Content.Find.Execute($needle5, $true, $false, $false, $false, $false, $true, 1, $false, $replace5, 1) | Out-Null

$needle6 = "in 1 bitmap, waardoor je dan alleen nog per karakter het goede stukje bitmap uit moet knippen"
$d.Content.Find.Execute($needle6, $true, $false, $false, $false, $false, $true, 1, $false, $needle6, 1) | Out-Null

# ---------------------------------------------------------------------------
# Move the _GoBack bookmark from the end of paragraph 1 to right after
# "...zoals inverteren en croppen" in paragraph 2 (just before the period).
# ---------------------------------------------------------------------------

$d.Bookmarks("_GoBack").Delete()

$target = $d.Content
$target.Find.Execute("zoals inverteren en croppen", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($target.End, $target.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Host "done"
